$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C
$ws.Range("C1").Value = "spv"

# Updated values for column B (existing "OT" series) and new column C ("spv")
# Each entry: @(B-value, C-value) for rows 2..25 (A2..A25 = 1..24, unchanged)
$values = @(
    @(0.45, 0),
    @(0.40500000000000003, 0),
    @(0.39599999999999996, 0),
    @(0.38442857142857151, 0),
    @(0.46285714285714297, 0.1),
    @(0.59073603737142866, 0.2),
    @(0.6457142857142858, 0.3),
    @(0.63505519657142873, 0.5),
    @(0.62439610742857143, 0.7),
    @(0.63505519657142873, 0.9),
    @(0.63150216685714289, 0.95),
    @(0.65637337485714298, 1),
    @(0.63150216685714289, 0.92),
    @(0.62511749999999999, 0.9),
    @(0.66535714285714287, 0.85),
    @(0.68889168707142856, 0.75),
    @(0.81861382350000012, 0.5),
    @(0.8566235515928573, 0.3),
    @(0.89463327968571438, 0.2),
    @(1, 0.1),
    @(0.8566235515928573, 0),
    @(0.77142857142857146, 0),
    @(0.60998541445714294, 0),
    @(0.53425654320000004, 0)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $entry = $values[$i]
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
}

# Update the active selection to match the edited workbook
[void]$ws.Range("Q20").Select()
